$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '30.166.52'
Set-TextValue 'E2' '  -0.66%  '
Set-TextValue 'D3' '1.905.68'
Set-TextValue 'E3' '  -1.59%  '
Set-TextValue 'E4' '  +0.29%  '
Set-TextValue 'D5' '0.7289'
Set-TextValue 'E5' '  -6.02%  '
Set-TextValue 'D6' '242.67'
Set-TextValue 'E6' '  -1.80%  '
Set-TextValue 'D7' '1.002'
Set-TextValue 'E7' '  +0.26%  '
Set-TextValue 'D8' '0.3117'
Set-TextValue 'E8' '  -3.08%  '
Set-TextValue 'D9' '26.30'
Set-TextValue 'E9' '  -6.09%  '
Set-TextValue 'D10' '0.06883'
Set-TextValue 'E10' '  -3.06%  '
Set-TextValue 'D11' '0.7748'
Set-TextValue 'E11' '  -1.21%  '
Set-TextValue 'D12' '0.07942'
Set-TextValue 'E12' '  -1.03%  '
Set-TextValue 'D13' '1.877.33'
Set-TextValue 'E13' '  -3.09%  '
Set-TextValue 'D14' '5.253'
Set-TextValue 'E14' '  -2.45%  '
Set-TextValue 'D15' '91.01'
Set-TextValue 'E15' '  -4.20%  '
Set-TextValue 'D16' '30.175.48'
Set-TextValue 'E16' '  -0.66%  '
Set-TextValue 'D17' '14.11'
Set-TextValue 'E17' '  -3.12%  '
Set-TextValue 'D18' '5.821'
Set-TextValue 'E18' '  -0.31%  '
Set-TextValue 'B19' 'ShibaInu'
Set-TextValue 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D19' '0.000007745'
Set-TextValue 'E19' '  -3.46%  '
Set-TextValue 'B20' 'BitcoinCash'
Set-TextValue 'C20' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D20' '238.51'
Set-TextValue 'E20' '  -6.72%  '
Set-TextValue 'D21' '1.001'
Set-TextValue 'E21' '  +0.22%  '
Set-TextValue 'D22' '2.146.90'
Set-TextValue 'E22' '  -2.07%  '
Set-TextValue 'E23' '  +0.26%  '
Set-TextValue 'D24' '7.001'
Set-TextValue 'E24' '  +3.59%  '
Set-TextValue 'D25' '9.347'
Set-TextValue 'E25' '  -2.78%  '
Set-TextValue 'E26' '  +0.40%  '
Set-TextValue 'D27' '19.00'
Set-TextValue 'E27' '  -0.71%  '
Set-TextValue 'D28' '0.1274'
Set-TextValue 'E28' '  -5.46%  '
Set-TextValue 'E29' '  -10.73%  '
Set-TextValue 'D30' '1.351'
Set-TextValue 'E30' '  -1.02%  '
Set-TextValue 'E31' '  +1.08%  '
Set-TextValue 'D32' '4.290'
Set-TextValue 'E32' '  -3.27%  '
Set-TextValue 'D33' '4.072'
Set-TextValue 'E33' '  -1.85%  '
Set-TextValue 'D34' '0.05114'
Set-TextValue 'E34' '  -1.65%  '
Set-TextValue 'D35' '1.285'
Set-TextValue 'E35' '  +0.00%  '
Set-TextValue 'D36' '0.7367'
Set-TextValue 'E36' '  -2.21%  '
Set-TextValue 'D37' '2.755'
Set-TextValue 'E37' '  -0.61%  '
Set-TextValue 'D38' '0.01925'
Set-TextValue 'E38' '  -2.67%  '
Set-TextValue 'D39' '2.782'
Set-TextValue 'E39' '  -1.09%  '
Set-TextValue 'D40' '6.342'
Set-TextValue 'E40' '  -2.34%  '
Set-TextValue 'D41' '74.51'
Set-TextValue 'E41' '  -5.90%  '
Set-TextValue 'D42' '0.4421'
Set-TextValue 'E42' '  -2.36%  '
Set-TextValue 'D43' '1.931'
Set-TextValue 'E43' '  -2.72%  '
Set-TextValue 'E44' '  +0.08%  '
Set-TextValue 'D45' '0.8337'
Set-TextValue 'E45' '  -0.21%  '
Set-TextValue 'D46' '101.02'
Set-TextValue 'E46' '  -0.36%  '
Set-TextValue 'D47' '7.555'
Set-TextValue 'E47' '  +0.75%  '
Set-TextValue 'D48' '9.702'
Set-TextValue 'E48' '  -1.50%  '
Set-TextValue 'D49' '37.55'
Set-TextValue 'E49' '  +0.10%  '
Set-TextValue 'D50' '939.91'
Set-TextValue 'E50' '  -4.33%  '
Set-TextValue 'D51' '0.1181'
Set-TextValue 'E51' '  -0.54%  '
